$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "ecosystems"

# Update row 2: PLA virgin
$ws.Range("A2").Value = "PLA_virgin"
$ws.Range("B2").Value = 0.577

# Update row 3: PLA recycled
$ws.Range("A3").Value = "PLA_recycled"
$ws.Range("B3").Value = 0.444

# Update row 4: PLA recycled industrial
$ws.Range("A4").Value = "PLA_recycled_industrial"
$ws.Range("B4").Value = 0.5600000000000001
